$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (existing rows 10-71 shift down to 11-72)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record
$ws.Cells.Item(10, 1).Value2  = 9
$ws.Cells.Item(10, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value2  = "Metropolitana"
$ws.Cells.Item(10, 4).Value2  = 44831
$ws.Cells.Item(10, 5).Value2  = 13
$ws.Cells.Item(10, 6).Value2  = 100112035
$ws.Cells.Item(10, 7).Value2  = "Bruselas (repollito)"
$ws.Cells.Item(10, 8).Value2  = "Sin especificar"
$ws.Cells.Item(10, 9).Value2  = "Primera"
$ws.Cells.Item(10, 10).Value2 = 20
$ws.Cells.Item(10, 11).Value2 = 19000
$ws.Cells.Item(10, 12).Value2 = 20000
$ws.Cells.Item(10, 13).Value2 = 19500
$ws.Cells.Item(10, 14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item(10, 15).Value2 = "Hijuelas"
$ws.Cells.Item(10, 16).Value2 = 1300
$ws.Cells.Item(10, 17).Value2 = 15
$ws.Cells.Item(10, 18).Value2 = "Hortaliza"

# Make the date cell in the new row match the date formatting of the rest of column D
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat

# Append a new weekly record as row 73 (right after the previous last row, 72)
$ws.Cells.Item(73, 1).Value2  = 9
$ws.Cells.Item(73, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(73, 3).Value2  = "Metropolitana"
$ws.Cells.Item(73, 4).Value2  = 44832
$ws.Cells.Item(73, 5).Value2  = 13
$ws.Cells.Item(73, 6).Value2  = 100112035
$ws.Cells.Item(73, 7).Value2  = "Bruselas (repollito)"
$ws.Cells.Item(73, 8).Value2  = "Sin especificar"
$ws.Cells.Item(73, 9).Value2  = "Primera"
$ws.Cells.Item(73, 10).Value2 = 22
$ws.Cells.Item(73, 11).Value2 = 20000
$ws.Cells.Item(73, 12).Value2 = 20000
$ws.Cells.Item(73, 13).Value2 = 20000
$ws.Cells.Item(73, 14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item(73, 15).Value2 = "Hijuelas"
$ws.Cells.Item(73, 16).Value2 = 1333
$ws.Cells.Item(73, 17).Value2 = 15
$ws.Cells.Item(73, 18).Value2 = "Hortaliza"

# Match the date formatting of the rest of column D for the appended row too
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
